$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "canonical SMILES"
$ws.Range("D3").Value = "c1ccc2c(c1)[nH]c[n+]2c3ccc(cc3)[O-]"
$ws.Range("D4").Value = "c1ccc2c(c1)ncn2c3ccc(cc3)O"
$ws.Range("D5").Value = "c1ccc2c(c1)[nH+]cn2c3ccc(cc3)O"
$ws.Range("D6").Value = "c1ccc2c(c1)ncn2c3ccc(cc3)[O-]"

$ws.Columns.Item(4).ColumnWidth = 35.92
